$wb = $excel.ActiveWorkbook

# Worksheets in this workbook:
#   1 -> "Overview"
#   2 -> "zh-cn"
#   3 -> "de-de"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Overview" sheet: Latest HO Xliff Generate Date column (G) for both data rows
# was regenerated a few minutes later.
$wsOverview.Range("G2").Value = "2016-09-07 15:13:46"
$wsOverview.Range("G3").Value = "2016-09-07 15:13:46"

# "zh-cn" sheet: Priority changed from human-translation ("ht") to machine
# translation ("mt"), and the handoff/handback timestamps moved forward.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-07 15:13:31"
$wsZhCn.Range("H3").Value = "2016-09-07 15:13:31"
$wsZhCn.Range("K2").Value = "2016-09-07 15:14:45"
$wsZhCn.Range("K3").Value = "2016-09-07 15:14:45"

# "de-de" sheet: same Priority and Latest HO Xliff Generate Date updates as
# above, plus a refreshed Correspond Handback DateTime.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-07 15:13:46"
$wsDeDe.Range("H3").Value = "2016-09-07 15:13:46"
$wsDeDe.Range("K2").Value = "2016-09-07 15:15:11"
$wsDeDe.Range("K3").Value = "2016-09-07 15:15:11"
